# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp
# - Update España (Spain) totals
# - Re-rank Indonesia / Bielorrusia / Noruega (Indonesia moves up with fresh numbers,
#   Bielorrusia and Noruega keep their previous numbers but shift down one row)
# - Small bump for Singapur, Malasia, Finlandia, Estonia
# - Re-rank Cuba / Afganistan (Afganistan moves up with fresh numbers, Cuba keeps its
#   previous numbers but shifts down one row)
# - Re-rank Vietnam / Mali / Tanzania (Tanzania moves up with fresh numbers, Vietnam and
#   Mali keep their previous numbers but shift down one row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 22 de Abril de 2020 a las 11:22"

# --- Row 5: Espana ------------------------------------------------------
$ws.Cells.Item(5, 2).Value = 208389
$ws.Cells.Item(5, 3).Value = 4211
$ws.Cells.Item(5, 4).Value = 85915
$ws.Cells.Item(5, 5).Value = 100757
$ws.Cells.Item(5, 6).Value = 7705
$ws.Cells.Item(5, 7).Value = 435
$ws.Cells.Item(5, 8).Value = 21717

# --- Row 31: Singapur ---------------------------------------------------
$ws.Cells.Item(31, 6).Value = 27

# --- Rows 38-40: Indonesia / Bielorrusia / Noruega reshuffle -----------
$ws.Cells.Item(38, 1).Value = "Indonesia"
$ws.Cells.Item(38, 2).Value = 7418
$ws.Cells.Item(38, 3).Value = 283
$ws.Cells.Item(38, 4).Value = 913
$ws.Cells.Item(38, 5).Value = 5870
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 19
$ws.Cells.Item(38, 8).Value = 635

$ws.Cells.Item(39, 1).Value = "Bielorrusia"
$ws.Cells.Item(39, 2).Value = 7281
$ws.Cells.Item(39, 3).Value = 558
$ws.Cells.Item(39, 4).Value = 769
$ws.Cells.Item(39, 5).Value = 6454
$ws.Cells.Item(39, 6).Value = 92
$ws.Cells.Item(39, 7).Value = 3
$ws.Cells.Item(39, 8).Value = 58

$ws.Cells.Item(40, 1).Value = "Noruega"
$ws.Cells.Item(40, 2).Value = 7241
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 32
$ws.Cells.Item(40, 5).Value = 7027
$ws.Cells.Item(40, 6).Value = 58
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 182

# --- Row 47: Malasia -----------------------------------------------------
$ws.Cells.Item(47, 2).Value = 5532
$ws.Cells.Item(47, 3).Value = 50
$ws.Cells.Item(47, 4).Value = 3452
$ws.Cells.Item(47, 5).Value = 1987
$ws.Cells.Item(47, 7).Value = 1
$ws.Cells.Item(47, 8).Value = 93

# --- Row 51: Finlandia ---------------------------------------------------
$ws.Cells.Item(51, 2).Value = 4129
$ws.Cells.Item(51, 3).Value = 115
$ws.Cells.Item(51, 5).Value = 1988

# --- Row 71: Estonia ------------------------------------------------------
$ws.Cells.Item(71, 2).Value = 1559
$ws.Cells.Item(71, 3).Value = 7
$ws.Cells.Item(71, 4).Value = 184
$ws.Cells.Item(71, 5).Value = 1331
$ws.Cells.Item(71, 6).Value = 7
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 44

# --- Rows 81-82: Cuba / Afganistan reshuffle -----------------------------
$ws.Cells.Item(81, 1).Value = "Afganistan"
$ws.Cells.Item(81, 2).Value = 1143
$ws.Cells.Item(81, 3).Value = 51
$ws.Cells.Item(81, 4).Value = 166
$ws.Cells.Item(81, 5).Value = 937
$ws.Cells.Item(81, 6).Value = 7
$ws.Cells.Item(81, 7).Value = 4
$ws.Cells.Item(81, 8).Value = 40

$ws.Cells.Item(82, 1).Value = "Cuba"
$ws.Cells.Item(82, 2).Value = 1137
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(82, 4).Value = 309
$ws.Cells.Item(82, 5).Value = 790
$ws.Cells.Item(82, 6).Value = 18
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 38

# --- Rows 122-124: Vietnam / Mali / Tanzania reshuffle -------------------
$ws.Cells.Item(122, 1).Value = "Tanzania"
$ws.Cells.Item(122, 2).Value = 284
$ws.Cells.Item(122, 3).Value = 30
$ws.Cells.Item(122, 4).Value = 11
$ws.Cells.Item(122, 5).Value = 263
$ws.Cells.Item(122, 6).Value = 7
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 10

$ws.Cells.Item(123, 1).Value = "Vietnam"
$ws.Cells.Item(123, 2).Value = 268
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 222
$ws.Cells.Item(123, 5).Value = 46
$ws.Cells.Item(123, 6).Value = 8
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 0

$ws.Cells.Item(124, 1).Value = "Mali"
$ws.Cells.Item(124, 2).Value = 258
$ws.Cells.Item(124, 3).Value = 0
$ws.Cells.Item(124, 4).Value = 57
$ws.Cells.Item(124, 5).Value = 187
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 14
